$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.779.74'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '1.700.10'
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("D4").Formula = "'1.003"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Formula = "'318.02"
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").Formula = "'1.004"
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Formula = "'0.3934"
$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("D8").Formula = "'0.4045"
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Formula = "'1.508"
$ws.Range("E9").Value = '  -2.96%  '

$ws.Range("D10").Formula = "'54.10"
$ws.Range("E10").Value = '  -2.23%  '

$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("D12").Formula = "'0.08911"
$ws.Range("E12").Value = '  +1.17%  '

$ws.Range("D13").Formula = "'7.263"
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").Formula = "'23.45"
$ws.Range("E14").Value = '  +0.39%  '

$ws.Range("D15").Formula = "'8.025"
$ws.Range("E15").Value = '  +4.41%  '

$ws.Range("D16").Formula = "'0.00001325"
$ws.Range("E16").Value = '  -0.57%  '

$ws.Range("D17").Value = '1.702.38'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").Formula = "'100.18"
$ws.Range("E18").Value = '  -1.44%  '

$ws.Range("D19").Formula = "'0.07035"
$ws.Range("E19").Value = '  -0.05%  '

$ws.Range("D20").Formula = "'19.65"
$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").Formula = "'7.007"
$ws.Range("E21").Value = '  +1.48%  '

$ws.Range("D22").Formula = "'1.002"
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("E23").Value = '  +2.49%  '

$ws.Range("D24").Value = '24.762.78'
$ws.Range("E24").Value = '  +0.72%  '

$ws.Range("D25").Formula = "'3.225"
$ws.Range("E25").Value = '  +8.58%  '

$ws.Range("D26").Formula = "'2.357"
$ws.Range("E26").Value = '  +0.63%  '

$ws.Range("D27").Formula = "'22.81"
$ws.Range("E27").Value = '  +1.70%  '

$ws.Range("D28").Formula = "'161.38"
$ws.Range("E28").Value = '  +0.50%  '

$ws.Range("D29").Formula = "'136.62"
$ws.Range("E29").Value = '  +2.06%  '

$ws.Range("D30").Formula = "'5.167"
$ws.Range("E30").Value = '  -1.43%  '

$ws.Range("D31").Formula = "'7.826"
$ws.Range("E31").Value = '  +3.05%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Formula = "'0.08743"
$ws.Range("E32").Value = '  +2.15%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Formula = "'1.080"
$ws.Range("E33").Value = '  -3.10%  '

$ws.Range("D34").Formula = "'7.180"
$ws.Range("E34").Value = '  -5.19%  '

$ws.Range("D35").Formula = "'11.27"
$ws.Range("E35").Value = '  +0.68%  '

$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").Formula = "'0.2754"
$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("B37").Value = 'WEMIXTOKEN'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Formula = "'1.952"
$ws.Range("E37").Value = '  -1.60%  '

$ws.Range("E38").Value = '  -2.70%  '

$ws.Range("D39").Formula = "'0.09197"
$ws.Range("E39").Value = '  +1.55%  '

$ws.Range("D40").Formula = "'0.02744"
$ws.Range("E40").Value = '  -1.64%  '

$ws.Range("D41").Formula = "'1.467"
$ws.Range("E41").Value = '  -0.14%  '

$ws.Range("D42").Formula = "'0.7705"
$ws.Range("E42").Value = '  -0.65%  '

$ws.Range("D43").Formula = "'15.87"
$ws.Range("E43").Value = '  +1.06%  '

$ws.Range("D44").Formula = "'0.7191"
$ws.Range("E44").Value = '  -1.29%  '

$ws.Range("D45").Formula = "'2.570"
$ws.Range("E45").Value = '  +2.35%  '

$ws.Range("D46").Formula = "'4.224"
$ws.Range("E46").Value = '  +0.81%  '

$ws.Range("D47").Formula = "'1.003"
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("D48").Formula = "'140.67"
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").Formula = "'1.315"
$ws.Range("E49").Value = '  +0.99%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Formula = "'0.07989"
$ws.Range("E50").Value = '  -0.29%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Formula = "'90.52"
$ws.Range("E51").Value = '  +2.42%  '
